$wb = $excel.ActiveWorkbook

# Work on the "Users list - M2M" sheet, which holds the test data rows.
$ws = $wb.Worksheets.Item("Users list - M2M")

# Row 3 currently references an id (2) that exists; change it to an id
# that does not exist, so the import is expected to fail.
$ws.Range("A3").Value = "ThatIdDoesntExist"

# Keep B3's text the same (it stays "Mitchell Admin Updated"); just
# touch it so the value is re-written through the shared string table.
$ws.Range("B3").Value = "Mitchell Admin Updated"

# Move the selection down to A2, matching the saved selection state.
[void]$ws.Range("A2").Select()

$ws.Activate()
